$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the left ("Javakheti" / Poti) table one more year: add a 2023
# column (K) that mirrors the formatting of the existing 2022 column (J).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1570
$ws.Range("K5").Value = 1204.5999999999999
$ws.Range("K6").Value = 1711.6
